# Apply updated symbol list values (Price / Volume(1h)) per the Mon Jan 23 20:24:39 UTC 2023 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "305.76" }
    @{ Cell = "E2"; Value = "0.55%" }
    @{ Cell = "D3"; Value = "36.25" }
    @{ Cell = "E3"; Value = "-1.48%" }
    @{ Cell = "D4"; Value = "5.042" }
    @{ Cell = "E4"; Value = "0.52%" }
    @{ Cell = "D5"; Value = "0.07880" }
    @{ Cell = "E5"; Value = "1.28%" }
    @{ Cell = "D6"; Value = "2.272" }
    @{ Cell = "E6"; Value = "7.70%" }
    @{ Cell = "D7"; Value = "7.994" }
    @{ Cell = "E7"; Value = "-0.36%" }
    @{ Cell = "D8"; Value = "0.9268" }
    @{ Cell = "E8"; Value = "0.69%" }
    @{ Cell = "D9"; Value = "0.09805" }
    @{ Cell = "E9"; Value = "-1.06%" }
    @{ Cell = "D10"; Value = "0.1875" }
    @{ Cell = "E10"; Value = "0.62%" }
    @{ Cell = "D11"; Value = "0.08993" }
    @{ Cell = "E11"; Value = "3.79%" }
    @{ Cell = "D12"; Value = "0.03742" }
    @{ Cell = "E12"; Value = "4.37%" }
    @{ Cell = "E13"; Value = "-0.70%" }
    @{ Cell = "D14"; Value = "0.001444" }
    @{ Cell = "E14"; Value = "-3.22%" }
    @{ Cell = "D15"; Value = "0.005614" }
    @{ Cell = "E15"; Value = "-1.39%" }
    @{ Cell = "D16"; Value = "3.460" }
    @{ Cell = "E16"; Value = "-0.02%" }
    @{ Cell = "D17"; Value = "4.144" }
    @{ Cell = "E17"; Value = "2.26%" }
    @{ Cell = "E18"; Value = "9.56%" }
    @{ Cell = "D19"; Value = "0.3367" }
    @{ Cell = "E19"; Value = "-2.19%" }
    @{ Cell = "D20"; Value = "0.1319" }
    @{ Cell = "E20"; Value = "1.08%" }
    @{ Cell = "D21"; Value = "5.132" }
    @{ Cell = "E21"; Value = "3.92%" }
    @{ Cell = "D22"; Value = "0.2251" }
    @{ Cell = "E22"; Value = "1.62%" }
    @{ Cell = "D23"; Value = "0.04581" }
    @{ Cell = "E23"; Value = "-0.70%" }
    @{ Cell = "E24"; Value = "-0.25%" }
    @{ Cell = "D25"; Value = "0.004774" }
    @{ Cell = "E25"; Value = "-7.18%" }
    @{ Cell = "D26"; Value = "0.0001302" }
    @{ Cell = "E26"; Value = "-7.44%" }
    @{ Cell = "E27"; Value = "73.65%" }
    @{ Cell = "D39"; Value = "0.01913" }
    @{ Cell = "E39"; Value = "6.68%" }
    @{ Cell = "D40"; Value = "0.04964" }
    @{ Cell = "E40"; Value = "6.01%" }
    @{ Cell = "D41"; Value = "0.007806" }
    @{ Cell = "E41"; Value = "1.30%" }
    @{ Cell = "D42"; Value = "0.1391" }
    @{ Cell = "E42"; Value = "-0.58%" }
    @{ Cell = "D43"; Value = "0.007813" }
    @{ Cell = "E43"; Value = "2.36%" }
    @{ Cell = "D44"; Value = "0.002144" }
    @{ Cell = "E44"; Value = "0.11%" }
    @{ Cell = "D45"; Value = "0.01143" }
    @{ Cell = "E45"; Value = "9.96%" }
    @{ Cell = "D46"; Value = "0.00006157" }
    @{ Cell = "E46"; Value = "-2.79%" }
    @{ Cell = "D47"; Value = "0.00000000751" }
    @{ Cell = "E47"; Value = "-0.36%" }
    @{ Cell = "D48"; Value = "51.77" }
    @{ Cell = "E48"; Value = "54.85%" }
    @{ Cell = "D49"; Value = "0.001802" }
    @{ Cell = "E49"; Value = "-10.37%" }
    @{ Cell = "D50"; Value = "0.00002104" }
    @{ Cell = "E50"; Value = "-0.36%" }
    @{ Cell = "D51"; Value = "0.0002004" }
    @{ Cell = "E51"; Value = "-0.36%" }
)

foreach ($u in $updates) {
    # Leading apostrophe forces text storage so numeric-looking / percent
    # strings (e.g. '305.76', '0.55%') stay strings, matching the source file
    # which stores every Price/Volume cell as inline text, not a number.
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
